$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores price strings as text in the source data (e.g. "247.15",
# "0.9989") rather than numbers. Mark each cell being updated as Text first so
# Excel keeps storing the new price as text instead of auto-converting it to a
# number, matching the workbooks existing layout.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.951.10'
$ws.Range("E2").Value = '  +8.11%  '
$ws.Range("D3").Value = '1.826.10'
$ws.Range("E3").Value = '  +5.77%  '
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '247.15'
$ws.Range("E5").Value = '  +2.72%  '
$ws.Range("D6").Value = '0.9986'
$ws.Range("D7").Value = '0.4945'
$ws.Range("E7").Value = '  +2.80%  '
$ws.Range("D8").Value = '43.76'
$ws.Range("E8").Value = '  +5.99%  '
$ws.Range("D9").Value = '0.2790'
$ws.Range("E9").Value = '  +7.65%  '
$ws.Range("D10").Value = '0.06411'
$ws.Range("E10").Value = '  +3.81%  '
$ws.Range("D11").Value = '1.816.73'
$ws.Range("E11").Value = '  +5.29%  '
$ws.Range("D12").Value = '16.79'
$ws.Range("E12").Value = '  +5.91%  '
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.6455'
$ws.Range("E14").Value = '  +7.00%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '84.25'
$ws.Range("E15").Value = '  +9.57%  '
$ws.Range("D16").Value = '4.681'
$ws.Range("E16").Value = '  +4.91%  '
$ws.Range("D17").Value = '28.974.69'
$ws.Range("E17").Value = '  +9.00%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").Value = '0.000007323'
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = '0.9982'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '12.24'
$ws.Range("E21").Value = '  +7.76%  '
$ws.Range("D22").Value = '2.041.19'
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("D23").Value = '4.570'
$ws.Range("E23").Value = '  +3.67%  '
$ws.Range("D24").Value = '8.863'
$ws.Range("E24").Value = '  +4.53%  '
$ws.Range("D25").Value = '5.369'
$ws.Range("E25").Value = '  +6.09%  '
$ws.Range("D26").Value = '144.59'
$ws.Range("E26").Value = '  +3.34%  '
$ws.Range("D27").Value = '129.46'
$ws.Range("E27").Value = '  +21.66%  '
$ws.Range("D28").Value = '16.43'
$ws.Range("E28").Value = '  +8.06%  '
$ws.Range("D29").Value = '1.888'
$ws.Range("E29").Value = '  +5.92%  '
$ws.Range("D30").Value = '1.403'
$ws.Range("E30").Value = '  +2.58%  '
$ws.Range("D31").Value = '4.137'
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("D32").Value = '0.08368'
$ws.Range("E32").Value = '  +5.55%  '
$ws.Range("D33").Value = '3.792'
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").Value = '0.04933'
$ws.Range("E34").Value = '  +9.01%  '
$ws.Range("D35").Value = '1.100'
$ws.Range("E35").Value = '  +10.02%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.704'
$ws.Range("E36").Value = '  +4.37%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6741'
$ws.Range("E37").Value = '  +9.04%  '
$ws.Range("D38").Value = '2.302'
$ws.Range("E38").Value = '  +15.42%  '
$ws.Range("D39").Value = '2.714'
$ws.Range("E39").Value = '  +11.12%  '
$ws.Range("D40").Value = '0.9478'
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("D41").Value = '6.152'
$ws.Range("E41").Value = '  +9.61%  '
$ws.Range("D42").Value = '0.01587'
$ws.Range("E42").Value = '  +6.22%  '
$ws.Range("D43").Value = '0.9990'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = '100.55'
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("D45").Value = '0.4086'
$ws.Range("E45").Value = '  +6.73%  '
$ws.Range("D46").Value = '7.195'
$ws.Range("E46").Value = '  +6.24%  '
$ws.Range("D47").Value = '0.1224'
$ws.Range("E47").Value = '  +5.98%  '
$ws.Range("D48").Value = '0.05523'
$ws.Range("E48").Value = '  +3.16%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.148'
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '31.69'
$ws.Range("E50").Value = '  +5.35%  '
$ws.Range("D51").Value = '1.305'
$ws.Range("E51").Value = '  +4.87%  '
